$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 187, shifting the existing
# rows 187:211 down to 190:214 (values/styles move with the rows).
$ws.Rows("187:189").Insert()

# Populate the 3 newly inserted rows with the new weekly price entries
# (date 2021-09-10 = serial 44449) for Kiwi Hayward, Especial/Primera/Segunda,
# Terminal La Palmera de La Serena / Region de O'Higgins.

$newRows = @(
    @{ Row=187; L="Especial"; M=20; N=365000; O=370000; P=367500; S=817 },
    @{ Row=188; L="Primera";  M=20; N=305000; O=310000; P=307500; S=683 },
    @{ Row=189; L="Segunda";  M=20; N=265000; O=270000; P=267500; S=594 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 8
    $ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = 44449
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100101
    $ws.Cells.Item($row, 8).Value = "Berries"
    $ws.Cells.Item($row, 9).Value = 100101007
    $ws.Cells.Item($row, 10).Value = "Kiwi"
    $ws.Cells.Item($row, 11).Value = "Hayward"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "`$/bins (450 kilos)"
    $ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 450
}
